$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition): update F3, F5, F6
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 3230
$wsExhibition.Range("F5").Value = 27
$wsExhibition.Range("F6").Value = 139

# Sheet "全部类型" (All types): update F7, F9, F11
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 3230
$wsAll.Range("F9").Value = 27
$wsAll.Range("F11").Value = 139
